$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 1262
$ws.Range("L3").Value = 1274
$ws.Range("I4").Value = 1829
$ws.Range("K4").Value = 1744
$ws.Range("L4").Value = 353
$ws.Range("L5").Value = 86
$ws.Range("K6").Value = 9121
$ws.Range("L6").Value = 1261
$ws.Range("I7").Value = 26294
$ws.Range("K7").Value = 27536
$ws.Range("L7").Value = 4236

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L5").Value = 12
$ws.Range("L7").Value = 130
$ws.Range("L8").Value = 260
$ws.Range("L9").Value = 23
$ws.Range("L11").Value = 70
$ws.Range("L14").Value = 20
$ws.Range("L15").Value = 29
$ws.Range("L18").Value = 32
$ws.Range("K19").Value = 793
$ws.Range("L19").Value = 128
$ws.Range("L20").Value = 114
$ws.Range("L22").Value = 15
$ws.Range("I25").Value = 132
$ws.Range("K27").Value = 261
$ws.Range("K29").Value = 1520
$ws.Range("L29").Value = 205
$ws.Range("L36").Value = 71
$ws.Range("L40").Value = 7
$ws.Range("L42").Value = 133
$ws.Range("L46").Value = 11
$ws.Range("L48").Value = 62
$ws.Range("L50").Value = 26
$ws.Range("L51").Value = 56
$ws.Range("L52").Value = 84
$ws.Range("L53").Value = 53
$ws.Range("L54").Value = 94
$ws.Range("K63").Value = 78
$ws.Range("L63").Value = 14
$ws.Range("L64").Value = 30
$ws.Range("L67").Value = 157
$ws.Range("L72").Value = 19
$ws.Range("L73").Value = 35
$ws.Range("L78").Value = 65
$ws.Range("L79").Value = 116
$ws.Range("L83").Value = 97
$ws.Range("L84").Value = 47
$ws.Range("L85").Value = 216
$ws.Range("L93").Value = 23
$ws.Range("L95").Value = 66
$ws.Range("L96").Value = 37
$ws.Range("L99").Value = 64
$ws.Range("I101").Value = 26294
$ws.Range("K101").Value = 27536
$ws.Range("L101").Value = 4236

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("L2").Value = 7
$ws.Range("L7").Value = 20

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L5").Value = 1
$ws.Range("L7").Value = 37

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L6").Value = 39
$ws.Range("L7").Value = 130

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L3").Value = 19
$ws.Range("L6").Value = 22
$ws.Range("L7").Value = 70

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 58
$ws.Range("L6").Value = 45
$ws.Range("L7").Value = 216

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L3").Value = 23
$ws.Range("L7").Value = 84

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("L3").Value = 10
$ws.Range("L7").Value = 53

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 75
$ws.Range("L3").Value = 81
$ws.Range("L4").Value = 21
$ws.Range("L7").Value = 260

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L3").Value = 44
$ws.Range("L7").Value = 97

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L2").Value = 24
$ws.Range("L7").Value = 66

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L2").Value = 15
$ws.Range("L7").Value = 64

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 47
$ws.Range("L4").Value = 15
$ws.Range("L6").Value = 43
$ws.Range("L7").Value = 157

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L3").Value = 17
$ws.Range("L7").Value = 47

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L2").Value = 25
$ws.Range("L4").Value = 6
$ws.Range("L7").Value = 94

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L3").Value = 69
$ws.Range("K6").Value = 448
$ws.Range("L6").Value = 56
$ws.Range("K7").Value = 1520
$ws.Range("L7").Value = 205

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L2").Value = 10
$ws.Range("L7").Value = 62

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 38
$ws.Range("L3").Value = 43
$ws.Range("K4").Value = 32
$ws.Range("K7").Value = 793
$ws.Range("L7").Value = 128

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L4").Value = 11
$ws.Range("L7").Value = 133

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L4").Value = 8
$ws.Range("L7").Value = 65

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("L2").Value = 2
$ws.Range("L7").Value = 11

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L2").Value = 39
$ws.Range("L3").Value = 37
$ws.Range("L7").Value = 116

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("L3").Value = 6
$ws.Range("L7").Value = 30

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L6").Value = 37
$ws.Range("L7").Value = 114

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("L2").Value = 11
$ws.Range("L7").Value = 32

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L5").Value = 1
$ws.Range("L7").Value = 71

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("L2").Value = 8
$ws.Range("L7").Value = 23

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("I4").Value = 5
$ws.Range("I7").Value = 132

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L6").Value = 6
$ws.Range("L7").Value = 29

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("L3").Value = 9
$ws.Range("L7").Value = 26

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("L2").Value = 5
$ws.Range("L7").Value = 23

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L3").Value = 9
$ws.Range("L7").Value = 35

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("L6").Value = 7
$ws.Range("L7").Value = 12

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K4").Value = 33
$ws.Range("K7").Value = 261

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L4").Value = 4
$ws.Range("L7").Value = 56

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("L2").Value = 7
$ws.Range("L3").Value = 5
$ws.Range("L7").Value = 15

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("L3").Value = 4
$ws.Range("L7").Value = 19

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("L5").Value = 1
$ws.Range("L7").Value = 7
